# Update NanoSplits v2.0.0 template: refreshed field descriptions (comments),
# expanded dataset_type vocabulary list, and updated pav:createdOn timestamp.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("NanoSplits")
$wsTypes = $wb.Worksheets.Item("dataset_type")
$wsMeta = $wb.Worksheets.Item(".metadata")

# ---- 1. Refresh the header comments on the NanoSplits sheet (A1:G1) ----
$wsMain.Range("A1").Comment.Text("(Required) The unique identifier from HuBMAP or SenNet for the sample (such as a
block, section, or suspension) used to perform the assay. For instance, in an
RNAseq assay, the parent sample would be the suspension, while in imaging
assays, it would be the tissue section. If the assay is derived from multiple
parent samples, this field should contain a comma-separated list of identifiers.
Example: HBM386.ZGKG.235, HBM672.MKPK.442")
$wsMain.Range("B1").Comment.Text("A locally assigned identifier provided by the data provider for the dataset. It
is used to reference an external metadata record that may be maintained
independently, enabling traceability and supporting provenance tracking.
Example: Visium_9OLC_A4_S1")
$wsMain.Range("C1").Comment.Text("(Required) The DOI for the protocols.io page that details the assay or the
procedures used for sample procurement and preparation. For example, in the case
of an imaging assay, the protocol may start with tissue section staining and end
with the generation of an OME-TIFF file. The documented protocol should also
include any image processing steps involved in producing the final OME-TIFF.
Example: https://dx.doi.org/10.17504/protocols.io.eq2lyno9qvx9/v1")
$wsMain.Range("D1").Comment.Text("(Required) The specific type of dataset being produced. Example: RNAseq")
$wsMain.Range("E1").Comment.Text("(Required) The name of the file containing the ORCID IDs for all contributors to
this dataset. Example: ./contributors.csv")
$wsMain.Range("F1").Comment.Text("(Required) The top-level directory containing the raw and/or processed data. For
a single dataset upload, this might be represented as `".`", whereas for a data
upload containing multiple datasets, this would be the directory name for the
respective dataset. For example, if the data is within a directory named
`"TEST001-RK`", use the syntax `"./TEST001-RK`" for this field. If there are
multiple directory levels, use the format `"./TEST001-RK/Run1/Pass2`", where
`"Pass2`" is the subdirectory where the single dataset's data is stored. This is
an internal metadata field used solely for data ingestion. Example: ./TEST001-RK")
$wsMain.Range("G1").Comment.Text("(Required) The unique string identifier for the metadata specification version,
which is easily interpretable by computers for purposes of data validation and
processing. Example: 22bc762a-5020-419d-b170-24253ed9e8d9")

# ---- 2. Rewrite the dataset_type lookup list with the refreshed HRAVS vocabulary ----
$datasetTypes = @(
    ,@("HiFi-Slide", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000195")
    ,@("SNARE-seq2", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000264")
    ,@("COMET", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000458")
    ,@("Visium (no probes)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000302")
    ,@("DESI", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000204")
    ,@("Confocal", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000206")
    ,@("Stereo-seq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000385")
    ,@("Visium (with probes)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000303")
    ,@("Molecular Cartography", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000217")
    ,@("DBiT-seq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000428")
    ,@("Seq-Scope", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000390")
    ,@("CosMx Transcriptomics", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000218")
    ,@("CyCIF", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000200")
    ,@("Light Sheet", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000168")
    ,@("seqFISH", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000397")
    ,@("ATACseq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000309")
    ,@("CosMx Proteomics", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000435")
    ,@("Singular Genomics G4X", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000429")
    ,@("Visium HD", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000451")
    ,@("MERFISH", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000221")
    ,@("10X Multiome", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000215")
    ,@("4i", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000447")
    ,@("PhenoCycler", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000199")
    ,@("Second Harmonic Generation (SHG)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000208")
    ,@("Thick section Multiphoton MxIF", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000207")
    ,@("CyTOF", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000407")
    ,@("Olink", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000441")
    ,@("MIBI", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000172")
    ,@("Auto-fluorescence", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000205")
    ,@("FACS", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000440")
    ,@("Xenium", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000219")
    ,@("SIMS", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000202")
    ,@("Cell DIVE", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000159")
    ,@("CODEX", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000160")
    ,@("GeoMx (NGS)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000300")
    ,@("MUSIC", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000328")
    ,@("Pixel-seqV2", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000450")
    ,@("MALDI", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000201")
    ,@("2D Imaging Mass Cytometry", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296")
    ,@("Histology", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000197")
    ,@("Enhanced Stimulated Raman Spectroscopy (SRS)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000209")
    ,@("DART-FISH", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000396")
    ,@("Resolve", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000384")
    ,@("RNAseq", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000310")
    ,@("LC-MS", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000194")
    ,@("nanoSPLITS", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000312")
    ,@("GeoMx (nCounter)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000301")
    ,@("RNAseq (with probes)", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000311")
    ,@("MS Lipidomics", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000405")
    ,@("MPLEx", "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000448")
)

for ($i = 0; $i -lt $datasetTypes.Count; $i++) {
    $row = $i + 1
    $wsTypes.Cells.Item($row, 1).Value = $datasetTypes[$i][0]
    $wsTypes.Cells.Item($row, 2).Value = $datasetTypes[$i][1]
}

# ---- 3. Point the dataset_type data validation at the expanded range (A1:A50) ----
$wsMain.Range("D2").Validation.Formula1 = "'dataset_type'!`$A`$1:`$A`$50"

# ---- 4. Bump the pav:createdOn timestamp on the .metadata sheet ----
$wsMeta.Range("C2").Value = "2025-10-20T14:52:17-07:00"

